{"js": "// Revision after the first week of class.\n//\n// 1) Add a new changelog entry right before the existing\n//    \"2019-06-01: Draft completed...\" list item.\n// 2) Update several grade-bundle percentages / chapter counts.\n\nconst body = context.document.body;\n\n// --- 1) Insert new changelog list item -------------------------------\nconst paras = body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\nlet changelogAnchor = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text.indexOf(\"2019-06-01: Draft completed\") !== -1) {\n    changelogAnchor = paras.items[i];\n    break;\n  }\n}\nif (!changelogAnchor) {\n  throw new Error(\"Could not find the 2019-06-01 changelog paragraph\");\n}\n// Inserting \"Before\" this paragraph copies its paragraph properties\n// (style \"Compact\" + numPr numId 1001), matching the target markup.\nchangelogAnchor.insertParagraph(\n  \"2019-06-05: Minor tweaks and typo corrections\",\n  \"Before\"\n);\nawait context.sync();\n\n// --- 2) Simple text replacements in the grade-bundle checklists ------\nconst replacements = [\n  [\n    \"\u2610 ORION diagnostics for 5 chapters\",\n    \"\u2610 ORION diagnostics for 4 chapters\",\n  ],\n  [\n    \"\u2610 ORION diagnostics for 7 chapters\",\n    \"\u2610 ORION diagnostics for 6 chapters\",\n  ],\n  [\n    \"\u2610 60% on two the following: ORION proficiency, Mastery, iClicker\",\n    \"\u2610 50% on two the following: ORION proficiency, Mastery, iClicker\",\n  ],\n  [\n    \"\u2610 ORION diagnostics for 9 chapters\",\n    \"\u2610 ORION diagnostics for 8 chapters\",\n  ],\n  [\n    \"\u2610 70% on two of the following: ORION proficiency, Mastery, iClicker\",\n    \"\u2610 60% on two of the following: ORION proficiency, Mastery, iClicker\",\n  ],\n  [\n    \"\u2610 ORION diagnostics for 11 chapters\",\n    \"\u2610 ORION diagnostics for 10 chapters\",\n  ],\n  [\n    \"\u2610 80% on two of the following: ORION proficiency, Mastery, iClicker\",\n    \"\u2610 70% on two of the following: ORION proficiency, Mastery, iClicker\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Revision after the first week of class.\n#\n# 1) Add a new changelog entry right before the existing\n#    \"2019-06-01: Draft completed...\" list item.\n# 2) Update several grade-bundle percentages / chapter counts.\n\n$d = $word.ActiveDocument\n\n# --- 1) Insert new changelog list item -------------------------------\n$idx = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text -like \"*2019-06-01: Draft completed*\") {\n        $idx = $i\n        break\n    }\n}\nif ($idx -eq -1) {\n    throw \"Could not find the 2019-06-01 changelog paragraph\"\n}\n\n$target = $d.Paragraphs($idx)\n$r = $target.Range\n$r.Collapse(1)          # wdCollapseStart\n$r.InsertParagraphBefore()\n\n# The newly inserted (still empty) paragraph now occupies the original\n# index; it inherited the \"Compact\" style + numId 1001 list numbering\n# from the paragraph it was inserted in front of.\n$newPara = $d.Paragraphs($idx)\n$newPara.Range.Text = \"2019-06-05: Minor tweaks and typo corrections\"\n\n# --- 2) Simple text replacements in the grade-bundle checklists ------\n$pairs = @(\n    @(\"ORION diagnostics for 5 chapters\", \"ORION diagnostics for 4 chapters\"),\n    @(\"ORION diagnostics for 7 chapters\", \"ORION diagnostics for 6 chapters\"),\n    @(\"60% on two the following: ORION proficiency, Mastery, iClicker\", \"50% on two the following: ORION proficiency, Mastery, iClicker\"),\n    @(\"ORION diagnostics for 9 chapters\", \"ORION diagnostics for 8 chapters\"),\n    @(\"70% on two of the following: ORION proficiency, Mastery, iClicker\", \"60% on two of the following: ORION proficiency, Mastery, iClicker\"),\n    @(\"ORION diagnostics for 11 chapters\", \"ORION diagnostics for 10 chapters\"),\n    @(\"80% on two of the following: ORION proficiency, Mastery, iClicker\", \"70% on two of the following: ORION proficiency, Mastery, iClicker\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $fr = $d.Content\n    $found = $fr.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        throw \"Find text not found: $old\"\n    }\n}\n"}
